# Auto-generated: update price/profit figures (columns H-N) for scheduled market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6070.278
$arr[0,1] = 350.5
$arr[0,2] = 8930.166999999999
$arr[0,3] = 350.5
$arr[0,4] = 8930.166999999999
$arr[0,5] = 89.5
$arr[0,6] = -9810.166999999999
$ws.Range("H41:N41").Value = $arr

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3711.3823
$arr[0,1] = 2048.2942
$arr[0,2] = 5374.4707
$arr[0,3] = 2048.2942
$arr[0,4] = 5374.4707
$arr[0,5] = -925.2941999999998
$arr[0,6] = -7620.4707
$ws.Range("H86:N86").Value = $arr

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3711.3823
$arr[0,1] = 2048.2942
$arr[0,2] = 5374.4707
$arr[0,3] = 10241.471
$arr[0,4] = 26872.3535
$arr[0,5] = -4625.471
$arr[0,6] = -38104.3535
$ws.Range("H89:N89").Value = $arr

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 581
$arr[0,1] = 276.35715
$arr[0,2] = 1434
$arr[0,3] = 276.35715
$arr[0,4] = 1434
$arr[0,5] = 1643.64285
$arr[0,6] = -5274
$ws.Range("H107:N107").Value = $arr

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 109658.1
$arr[0,1] = 121286.78
$arr[0,2] = 5000
$arr[0,3] = 121286.78
$arr[0,4] = 5000
$arr[0,5] = -117844.78
$arr[0,6] = -11884
$ws.Range("H116:N116").Value = $arr

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 91548.91
$arr[0,1] = 100612.9
$arr[0,2] = 909
$arr[0,3] = 301838.7
$arr[0,4] = 2727
$arr[0,5] = -300181.7
$arr[0,6] = -6041
$ws.Range("H118:N118").Value = $arr

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 601.2143
$arr[0,1] = 354.54544
$arr[0,2] = 1505.6666
$arr[0,3] = 1063.63632
$arr[0,4] = 4516.9998
$arr[0,5] = 3936.36368
$arr[0,6] = -14516.9998
$ws.Range("H129:N129").Value = $arr

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2511.7585
$arr[0,1] = 1369.6
$arr[0,2] = 6465.385
$arr[0,3] = 4108.799999999999
$arr[0,4] = 19396.155
$arr[0,5] = -1578.799999999999
$arr[0,6] = -24456.155
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2586.1777
$arr[0,1] = 2575.3948
$arr[0,2] = 2644.7144
$arr[0,3] = 7726.1844
$arr[0,4] = 7934.1432
$arr[0,5] = -5176.1844
$arr[0,6] = -13034.1432
$ws.Range("H137:N137").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3599.5676
$arr[0,1] = 2929.2666
$arr[0,2] = 4056.5908
$arr[0,3] = 2929.2666
$arr[0,4] = 4056.5908
$arr[0,5] = -2717.2666
$arr[0,6] = -4480.5908
$ws.Range("H61:N61").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4750.1055
$arr[0,1] = 9571.429
$arr[0,2] = 1937.6666
$arr[0,3] = 9571.429
$arr[0,4] = 1937.6666
$arr[0,5] = -8885.429
$arr[0,6] = -3309.6666
$ws.Range("H63:N63").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4750.1055
$arr[0,1] = 9571.429
$arr[0,2] = 1937.6666
$arr[0,3] = 47857.145
$arr[0,4] = 9688.333000000001
$arr[0,5] = -44425.145
$arr[0,6] = -16552.333
$ws.Range("H66:N66").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2181.625
$arr[0,1] = 1263.9231
$arr[0,2] = 3885.9285
$arr[0,3] = 1263.9231
$arr[0,4] = 3885.9285
$arr[0,5] = -389.9231
$arr[0,6] = -5633.9285
$ws.Range("H74:N74").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2181.625
$arr[0,1] = 1263.9231
$arr[0,2] = 3885.9285
$arr[0,3] = 6319.6155
$arr[0,4] = 19429.6425
$arr[0,5] = -1951.6155
$arr[0,6] = -28165.6425
$ws.Range("H77:N77").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 22212.25
$arr[0,1] = 0
$arr[0,2] = 22212.25
$arr[0,3] = 0
$arr[0,4] = 22212.25
$arr[0,5] = $null
$arr[0,6] = -31888.25
$ws.Range("H119:N119").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2403.0952
$arr[0,1] = 2719.8
$arr[0,2] = 2115.182
$arr[0,3] = 8159.400000000001
$arr[0,4] = 6345.545999999999
$arr[0,5] = -5709.400000000001
$arr[0,6] = -11245.546
$ws.Range("H122:N122").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3107.6792
$arr[0,1] = 2944.7144
$arr[0,2] = 3424.5557
$arr[0,3] = 8834.143199999999
$arr[0,4] = 10273.6671
$arr[0,5] = -6304.143199999999
$arr[0,6] = -15333.6671
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3599.5676
$arr[0,1] = 2929.2666
$arr[0,2] = 4056.5908
$arr[0,3] = 8787.799800000001
$arr[0,4] = 12169.7724
$arr[0,5] = -6237.799800000001
$arr[0,6] = -17269.7724
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 25778
$arr[0,1] = 0
$arr[0,2] = 25778
$arr[0,3] = 0
$arr[0,4] = 25778
$arr[0,5] = $null
$arr[0,6] = -35898
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 37135
$arr[0,1] = 0
$arr[0,2] = 37135
$arr[0,3] = 0
$arr[0,4] = 37135
$arr[0,5] = $null
$arr[0,6] = -47415
$ws.Range("H139:N139").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2710.456
$arr[0,1] = 1659.0857
$arr[0,2] = 4383.091
$arr[0,3] = 1659.0857
$arr[0,4] = 4383.091
$arr[0,5] = -1364.0857
$arr[0,6] = -4973.091
$ws.Range("H31:N31").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2710.456
$arr[0,1] = 1659.0857
$arr[0,2] = 4383.091
$arr[0,3] = 1659.0857
$arr[0,4] = 4383.091
$arr[0,5] = -1457.0857
$arr[0,6] = -4787.091
$ws.Range("H34:N34").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1627.9143
$arr[0,1] = 918.15
$arr[0,2] = 2574.2666
$arr[0,3] = 918.15
$arr[0,4] = 2574.2666
$arr[0,5] = -715.15
$arr[0,6] = -2980.2666
$ws.Range("H58:N58").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 956.7778
$arr[0,1] = 960
$arr[0,2] = 952.75
$arr[0,3] = 2880
$arr[0,4] = 2858.25
$arr[0,5] = -430
$arr[0,6] = -7758.25
$ws.Range("H122:N122").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1627.9143
$arr[0,1] = 918.15
$arr[0,2] = 2574.2666
$arr[0,3] = 2754.45
$arr[0,4] = 7722.7998
$arr[0,5] = -204.4499999999998
$arr[0,6] = -12822.7998
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3000
$arr[0,1] = 0
$arr[0,2] = 3000
$arr[0,3] = 0
$arr[0,4] = 9000
$arr[0,5] = $null
$arr[0,6] = -9566
$ws.Range("H32:N32").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 463.95312
$arr[0,1] = 445.3611
$arr[0,2] = 487.85715
$arr[0,3] = 1336.0833
$arr[0,4] = 1463.57145
$arr[0,5] = 833.9167
$arr[0,6] = -5803.571449999999
$ws.Range("H113:N113").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3110.5
$arr[0,1] = 0
$arr[0,2] = 3110.5
$arr[0,3] = 0
$arr[0,4] = 9331.5
$arr[0,5] = $null
$arr[0,6] = -19251.5
$ws.Range("H127:N127").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1012.70514
$arr[0,1] = 3660
$arr[0,2] = 906.81335
$arr[0,3] = 10980
$arr[0,4] = 2720.44005
$arr[0,5] = -5940
$arr[0,6] = -12800.44005
$ws.Range("H131:N131").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1832.2222
$arr[0,1] = 865
$arr[0,2] = 3766.6667
$arr[0,3] = 2595
$arr[0,4] = 11300.0001
$arr[0,5] = 2545
$arr[0,6] = -21580.0001
$ws.Range("H139:N139").Value = $arr

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 37132.5
$arr[0,1] = 0
$arr[0,2] = 37132.5
$arr[0,3] = 0
$arr[0,4] = 37132.5
$arr[0,5] = $null
$arr[0,6] = -45312.5
$ws.Range("H110:N110").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2776
$arr[0,1] = 2782.5
$arr[0,2] = 2750
$arr[0,3] = 2782.5
$arr[0,4] = 2750
$arr[0,5] = -2580.5
$arr[0,6] = -3154
$ws.Range("H61:N61").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2776
$arr[0,1] = 2782.5
$arr[0,2] = 2750
$arr[0,3] = 2782.5
$arr[0,4] = 2750
$arr[0,5] = -612.5
$arr[0,6] = -7090
$ws.Range("H113:N113").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H116:N116").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4139.5884
$arr[0,1] = 2242.8965
$arr[0,2] = 6639.773
$arr[0,3] = 6728.689499999999
$arr[0,4] = 19919.319
$arr[0,5] = -4178.689499999999
$arr[0,6] = -25019.319
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H116:N116").Value = $arr

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2035.0555
$arr[0,1] = 905.6
$arr[0,2] = 2699.4412
$arr[0,3] = 2716.8
$arr[0,4] = 8098.323600000001
$arr[0,5] = -186.8000000000002
$arr[0,6] = -13158.3236
$ws.Range("H132:N132").Value = $arr
